$d = $word.ActiveDocument

# Change 2: add citation after "researches" in the classification sentence
$d.Content.Find.Execute(
    "Then we do classification by rules from few researches.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Then we do classification by rules from few researches(For example, http://www.pyoudeyer.com/emotionsIJHCS.pdf).",
    2) | Out-Null

# Change 3: loudness -> intensity (only for "Implement loudness normalizing")
$d.Content.Find.Execute(
    "Implement loudness normalizing",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implement intensity normalizing",
    2) | Out-Null

# Change 4: reword the speech-pauses bullet
$d.Content.Find.Execute(
    "Set optimal threshold to find speech pauses",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Find speech pauses and split speech to phrases",
    2) | Out-Null

Write-Output "done"
